# NpcShop.xlsx - add a MoneyType column (#42 npcshop support differenet types of resource cost)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the resource list for row 6 (Id 44010001): old test ids -> real shop ids
$ws.Range("B6").Value = "22031002;1|22031003;1|22031004;1|22031005;1"

# 2) Extend the "表1" table with a new MoneyType column (C)
$lo = $ws.ListObjects.Item(1)
$col = $lo.ListColumns.Add()
$col.Name = "MoneyType"

# 3) Header / type / description rows for column C, matching columns A & B styling
$ws.Range("C1").Value = "MoneyType"

$ws.Range("B2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("C2").Value = "int"

$ws.Range("B3").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C3").Value = "货币类型"

# 4) Data rows: MoneyType = 0 (default) except row 6 which now costs a special currency (1)
$ws.Range("C4").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("C6").Value = 1
$ws.Range("C7").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("C10").Value = 0

# 5) Print setup (A4 portrait) as set on the author's machine
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# 6) Selection cursor parked on C6, matching the saved view state
[void]$ws.Range("C6").Select()
